$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in / correct the "Absent" column (H) values to form the
# consolidated report.
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H16").Value = 0
